# feat: add 2022-Q4 data
#
# Before: sheet "总计" (summary) + sheet "2022-Q2" (detail for 2022-Q2).
# After : sheet "总计" (summary, now with a 2022-Q4 row added before the
#         2022-Q2 row) + sheet "2022-Q4" (new detail sheet, inserted right
#         after "总计") + sheet "2022-Q2" (unchanged detail sheet, moved
#         after "2022-Q4").

$wb = $excel.ActiveWorkbook

$summary = $wb.Worksheets.Item(1)      # "总计"
$q2      = $wb.Worksheets.Item("2022-Q2")

# ---------------------------------------------------------------------
# 1) Duplicate the existing "2022-Q2" detail sheet so the original data
#    survives untouched under the same name, right after itself. Then
#    rename the *original* sheet object to "2022-Q4" - it will be
#    overwritten with the new quarter's figures below.
# ---------------------------------------------------------------------
$q2.Copy($null, $q2)
$q2Dup = $wb.Worksheets.Item($wb.Worksheets.Count)
$q2Dup.Name = "2022-Q2-NEWCOPY"

$q4 = $q2
$q4.Name = "2022-Q4"

$q2Dup.Name = "2022-Q2"

# ---------------------------------------------------------------------
# 2) Rebuild the "2022-Q4" sheet contents (it currently still holds the
#    old 2022-Q2 rows because it is the renamed original sheet).
# ---------------------------------------------------------------------
$q4.Cells.Clear()

$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

$q4.Range("A2").Value = 0
$q4.Range("B2:G2").NumberFormat = "@"
$q4.Range("B2").Value = "486001"
$q4.Range("C2").Value = "工银瑞信中国机会全球配置股票（QDII）人民币"
$q4.Range("D2").Value = "6.42"
$q4.Range("E2").Value = "93.86"
$q4.Range("F2").Value = "1.90"
$q4.Range("G2").Value = "0.1220"
$q4.Range("H2").Value = 7

$q4.Range("A3").Value = 1
$q4.Range("B3:G3").NumberFormat = "@"
$q4.Range("B3").Value = "009562"
$q4.Range("C3").Value = "工银全球股票（QDII）美元"
$q4.Range("D3").Value = "6.42"
$q4.Range("E3").Value = "93.86"
$q4.Range("F3").Value = "1.90"
$q4.Range("G3").Value = "0.1220"
$q4.Range("H3").Value = 7

$q4.Range("A4").Value = 2
$q4.Range("B4:G4").NumberFormat = "@"
$q4.Range("B4").Value = "009563"
$q4.Range("C4").Value = "工银全球股票（QDII）港币"
$q4.Range("D4").Value = "6.42"
$q4.Range("E4").Value = "93.86"
$q4.Range("F4").Value = "1.90"
$q4.Range("G4").Value = "0.1220"
$q4.Range("H4").Value = 7

# Match the look & feel of the rest of the workbook: copy the bold /
# bordered header style from the "总计" sheet onto the new header row
# and first column (same style used for equivalent cells elsewhere).
$summary.Range("B1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)
$summary.Range("A2").Copy()
$q4.Range("A2:A4").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 3) Update the "总计" summary sheet: push the existing 2022-Q2 row down
#    to row 3 and put the new 2022-Q4 figures in row 2.
# ---------------------------------------------------------------------
$summary.Range("A2:D2").Copy($summary.Range("A3"))

$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 3
$summary.Range("D2").Value = 0.37

$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q2"
$summary.Range("C3").Value = 2
$summary.Range("D3").Value = 0.19
